# Add data for 2021-09-25
# Bumps the "Through 2021-09-16" snapshot to "Through 2021-09-17" (sheet
# name + the header cell describing the September-2021 column), and adds
# the newly-recorded carjacking counts for a handful of neighborhood / month
# combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet and update the "through" date label in B1 ---
$ws.Name = "Through 2021-09-17"
$ws.Range("B1").Value = "September 2021 (through September 17)"

# --- Updated counts (existing cells whose value increased by one) ---
$ws.Range("T3").Value = 4     # North Lawndale, September 2019
$ws.Range("AC4").Value = 2    # Humboldt Park, September 2018
$ws.Range("B5").Value = 7     # Austin, September 2021
$ws.Range("T17").Value = 3    # South Shore, September 2019
$ws.Range("K18").Value = 3    # Grand Boulevard, September 2020
$ws.Range("K31").Value = 5    # West Loop, September 2020

# --- Newly populated cells (previously empty, now a count of 1) ---
$ws.Range("B14").Value = 1    # Lower West Side, September 2021
$ws.Range("B21").Value = 1    # River North, September 2021
$ws.Range("AU30").Value = 1   # Lincoln Park, September 2016
$ws.Range("B33").Value = 1    # Lake View, September 2021
$ws.Range("B50").Value = 1    # Morgan Park, September 2021
$ws.Range("T54").Value = 1    # Hermosa, September 2019
$ws.Range("B56").Value = 1    # Edgewater, September 2021
$ws.Range("AC56").Value = 1   # Edgewater, September 2018
$ws.Range("T76").Value = 1    # Hegewisch, September 2019
